$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.866.06"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.512.03"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.39"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.88"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.511.34"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.00"
$ws.Range("E11").Value = "  +5.56%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.108.88"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.90"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.514.29"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.979.91"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("E19").Value = "  +8.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.45"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "438.18"
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.609"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.55"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.650.54"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.78"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.31"
$ws.Range("E29").Value = "  -4.28%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.63"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.55"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.97"
$ws.Range("E35").Value = "  -2.81%  "
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.04"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.36"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0896"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.06"
$ws.Range("E43").Value = "  -10.99%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.18"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.09"
$ws.Range("E46").Value = "  -7.76%  "
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.46"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.991"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("E51").Value = "  -1.26%  "
